# Edit PranshuGuptaResume.docx per commit "updated the new position details"
$d = $word.ActiveDocument

# -------------------------------------------------------------------
# Change 1: Azure OpenAI bullet — describe the data-plane work in detail
# "of components in the Azure OpenAI service." ->
#   "in the data plane of Azure OpenAI service, which handles all the
#    inference requests sent to OpenAI models hosted on Azure, via
#    Copilot and customers' AI agents."
# -------------------------------------------------------------------
$rng = $d.Content
$found = $rng.Find.Execute("of components in the Azure OpenAI service.")
if (-not $found) { throw "Change 1: target sentence not found" }

# Clear the old text, then insert each new run individually so they
# land as distinct <w:r> elements (mirrors the authored diff).
$rng.Text = ""
$rng.Collapse(1)

$rng.InsertAfter("in")
$rng.Collapse(0)

$rng.InsertAfter(" ")
$rng.Collapse(0)

$rng.InsertAfter("the data plane ")
$rng.Collapse(0)

$rng.InsertAfter("of Azure OpenAI service, ")
$rng.Collapse(0)

$rng.InsertAfter("which handles all the ")
$rng.Collapse(0)

$rng.InsertAfter("inference requests")
$rng.Collapse(0)

$rng.InsertAfter(" sen")
$rng.Collapse(0)

$rng.InsertAfter("t")
$rng.Collapse(0)

$rng.InsertAfter(" to OpenAI models hosted on Azure, via ")
$rng.Collapse(0)

$rng.InsertAfter("C")
$rng.Collapse(0)

$rng.InsertAfter("opilot")
$rng.Collapse(0)

$rng.InsertAfter(" and customers" + [char]0x2019 + " AI agents")
$rng.Collapse(0)

$rng.InsertAfter(".")
$rng.Collapse(0)

# -------------------------------------------------------------------
# Change 2: "multi-semester roadmap" bullet — trim the trailing detail
# ", by collecting data on the current state of the service,
#   identifying opportunities for improvements, and proposing solutions" ->
#   " control plane"
# (the closing "." run is left untouched)
# -------------------------------------------------------------------
$rng2 = $d.Content
$found2 = $rng2.Find.Execute(", by collecting data on the current state of the service, identifying opportunities for improvements, and proposing solutions")
if (-not $found2) { throw "Change 2: target sentence not found" }

$rng2.Text = ""
$rng2.Collapse(1)
$rng2.InsertAfter(" control plane")
$rng2.Collapse(0)

Write-Output $d.Content.Text
